# Apply the "less impactful hardening process" update to the
# Process Contribution - Characterization workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (Sheet7 -> Sheet4)
$ws.Name = "Sheet4"

# Update the report's generation date/time (header row 1)
$ws.Range("D1").Value = 45572
$ws.Range("F1").Value = 0.808306076388889

# Rename the "Induction Hardening Bearing Surface 1" process to
# "Induction Hardening Bearing Surfaces 1, 2" everywhere it appears:
#  - column header in the contribution table (row 17)
#  - process row label (row 21, column B)
$ws.Range("W17").Value = "Induction Hardening Bearing Surfaces 1, 2"
$ws.Range("B21").Value = "Induction Hardening Bearing Surfaces 1, 2"

# Reduce the impact contributed by the (now less impactful) hardening
# process, and update the overall total accordingly.
$ws.Range("E18").Value = 174.399516254188
$ws.Range("W18").Value = 34.074702596165
$ws.Range("E21").Value = 34.074702596165
$ws.Range("W21").Value = 34.074702596165
